$d = $word.ActiveDocument

# Locate the paragraph that begins the long Greek description starting
# with "Συμμετέχετε" (the GlobeAtNight activity description paragraph).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Συμμετέχετε")) {
        $target = $p
        break
    }
}

$r = $target.Range
# Exclude the trailing paragraph mark from the range so we only replace
# the paragraph's text content.
$r.End = $r.End - 1

$newText = "Συμμετέχετε σε μία παγκόσμια καμπάνια για να παρατηρήσετε και να καταγράψετε τη φωτεινότητα των πιο αμυδρά ορατών άστρων σαν μέσο για την μέτρηση της Φωτορρύπανσης σε μία δεδομένη περιοχή. Με τον εντοπισμό και την παρατήρηση του  Αστερισμός Λέων στον νυχτερινό ουρανό καθώς και με την σύγκριση των ανωτέρω με τα διαγράμματα για τα μεγέθη των άστρων,  άνθρωποι από όλον τον κόσμο θα μάθουν πώς τα φώτα στην κοινότητά τους συμβάλλουν στην Φωτορρύπανση. Με την κατάθεση των πορισμάτων τους στην ιστοσελίδα θα δημιουργηθεί ένα αρχείο σχετικά με το τι μπορεί να δει κανείς στον νυχτερινό ουρανό."

# Delete all the existing (many, differently-formatted) runs and insert a
# single new plain run containing the full updated text.
$r.Delete()
$r.InsertAfter($newText)
